$d = $word.ActiveDocument

# --- Change 1 & 2 -------------------------------------------------------
# In the "Post condiciones" value table (2nd table in the document), the
# row holding the post-condition text ("La postulación queda registrada
# en el sistema.") and its page-reference cell ("5") are emptied out
# while keeping the surrounding paragraph formatting intact.
$postCondTable = $d.Tables.Item(2)
$postCondRow = $postCondTable.Rows.Item(8)

$cellA = $postCondRow.Cells.Item(1).Range
$rangeA = $d.Range($cellA.Start, $cellA.End)
$rangeA.Text = ""

$cellB = $postCondRow.Cells.Item(2).Range
$rangeB = $d.Range($cellB.Start, $cellB.End)
$rangeB.Text = ""

# --- Change 3 -------------------------------------------------------
# Add a new precondition bullet right after "... está logeado a su
# cuenta." reusing the same list-paragraph formatting.
# NOTE: use $d.Content.Paragraphs (not the bare $d.Paragraphs collection)
# for paragraph lookups here - after touching $d.Tables the cached
# $d.Paragraphs collection can get anchored to the table and mis-report
# text for unrelated indices, while re-deriving the collection from
# $d.Content stays accurate.
$targetParaIndex = 0
$paras = $d.Content.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*está logeado a su cuenta.*") {
        $targetParaIndex = $i
        break
    }
}

$targetPara = $d.Content.Paragraphs.Item($targetParaIndex)
$targetPara.Range.InsertParagraphAfter()

$newPara = $d.Content.Paragraphs.Item($targetParaIndex + 1)
$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newRange.Text = "El usuario tiene los datos que actualizará o  agregará."
